$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells whose new value looks like a plain number need the cell
# number-format forced to Text first, otherwise Excel auto-converts the
# entry into a floating point number (losing the original text formatting).

$ws.Range('D2').Value = '65.512.03'
$ws.Range('E2').Value = '  +2.94%  '
$ws.Range('D3').Value = '3.404.60'
$ws.Range('E3').Value = '  +2.34%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '562.85'
$ws.Range('E5').Value = '  +3.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.19'
$ws.Range('E6').Value = '  +2.61%  '
$ws.Range('E7').Value = '  +3.02%  '
$ws.Range('D8').Value = '3.398.26'
$ws.Range('E8').Value = '  +2.42%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('E10').Value = '  +13.11%  '
$ws.Range('E11').Value = '  +3.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.94'
$ws.Range('E12').Value = '  +2.90%  '
$ws.Range('E13').Value = '  +6.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.17'
$ws.Range('E14').Value = '  +3.12%  '
$ws.Range('D15').Value = '3.950.23'
$ws.Range('E15').Value = '  +2.39%  '
$ws.Range('E16').Value = '  +3.26%  '
$ws.Range('D17').Value = '3.425.68'
$ws.Range('E17').Value = '  +3.20%  '
$ws.Range('E18').Value = '  +2.01%  '
$ws.Range('D19').Value = '65.562.65'
$ws.Range('E19').Value = '  +3.16%  '
$ws.Range('E20').Value = '  +2.29%  '
$ws.Range('E21').Value = '  +2.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '471.80'
$ws.Range('E22').Value = '  +14.70%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.09'
$ws.Range('E23').Value = '  +18.03%  '
$ws.Range('E24').Value = '  +3.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '86.94'
$ws.Range('E25').Value = '  +4.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.41'
$ws.Range('E26').Value = '  -2.20%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.94'
$ws.Range('E27').Value = '  +3.65%  '
$ws.Range('E28').Value = '  +6.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.89'
$ws.Range('E29').Value = '  +3.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.26'
$ws.Range('E30').Value = '  +7.88%  '
$ws.Range('E31').Value = '  +5.93%  '
$ws.Range('E32').Value = '  +2.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '62.72'
$ws.Range('E33').Value = '  +8.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '574.29'
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('E35').Value = '  +2.64%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  -4.31%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.53'
$ws.Range('E38').Value = '  +3.11%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.93'
$ws.Range('E39').Value = '  +2.55%  '
$ws.Range('D40').Value = '0.0₃0761'
$ws.Range('E40').Value = '  +3.50%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.374'
$ws.Range('E41').Value = '  +2.46%  '
$ws.Range('D42').Value = '3.090.74'
$ws.Range('E42').Value = '  -1.37%  '
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('E44').Value = '  +2.71%  '
$ws.Range('E45').Value = '  +4.62%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.50'
$ws.Range('E46').Value = '  +4.07%  '
$ws.Range('E47').Value = '  +6.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.21'
$ws.Range('E48').Value = '  -1.29%  '
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '136.90'
$ws.Range('E50').Value = '  +3.71%  '
$ws.Range('E51').Value = '  +4.32%  '
